$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Existing value updates (R/S columns) ---
$ws.Range("S2").Value = 2000

$ws.Range("R87").Value = 400
$ws.Range("S87").Value = 2200

$ws.Range("S88").Value = 2300

$ws.Range("S90").Value = 2400

$ws.Range("S93").Value = 2000

$ws.Range("R94").Value = 300
$ws.Range("S94").Value = 2100

$ws.Range("S98").Value = 2500

# --- New cells for rows 99-125: R (Obj-posXpix), S (Obj-posYpix), W (run), X (quality) ---

$ws.Range("R99").Value = 200
$ws.Range("S99").Value = 1800
$ws.Range("W99").Value = 1
$ws.Range("X99").Value = 1

$ws.Range("R100").Value = 200
$ws.Range("S100").Value = 1800
$ws.Range("W100").Value = 1
$ws.Range("X100").Value = 1

$ws.Range("R101").Value = 200
$ws.Range("S101").Value = 2100
$ws.Range("W101").Value = 1
$ws.Range("X101").Value = 1

$ws.Range("R102").Value = 200
$ws.Range("S102").Value = 2200
$ws.Range("W102").Value = 1
$ws.Range("X102").Value = 1

$ws.Range("R103").Value = 200
$ws.Range("S103").Value = 2200
$ws.Range("W103").Value = 1
$ws.Range("X103").Value = 1

$ws.Range("R104").Value = 200
$ws.Range("S104").Value = 2200
$ws.Range("W104").Value = 1
$ws.Range("X104").Value = 1

$ws.Range("R105").Value = 200
$ws.Range("S105").Value = 2500
$ws.Range("W105").Value = 1
$ws.Range("X105").Value = 1

$ws.Range("R106").Value = 100
$ws.Range("S106").Value = 1800
$ws.Range("W106").Value = 1
$ws.Range("X106").Value = 1

$ws.Range("R107").Value = 100
$ws.Range("S107").Value = 1900
$ws.Range("W107").Value = 1
$ws.Range("X107").Value = 1

$ws.Range("R108").Value = 100
$ws.Range("S108").Value = 2100
$ws.Range("W108").Value = 1
$ws.Range("X108").Value = 1

$ws.Range("R109").Value = 120
$ws.Range("S109").Value = 2100
$ws.Range("W109").Value = 1
$ws.Range("X109").Value = 1

$ws.Range("R110").Value = 200
$ws.Range("S110").Value = 2200
$ws.Range("W110").Value = 1
$ws.Range("X110").Value = 1

$ws.Range("R111").Value = 150
$ws.Range("S111").Value = 2300
$ws.Range("W111").Value = 1
$ws.Range("X111").Value = 1

$ws.Range("R112").Value = 150
$ws.Range("S112").Value = 2400
$ws.Range("W112").Value = 1
$ws.Range("X112").Value = 1

$ws.Range("R113").Value = 50
$ws.Range("S113").Value = 2300
$ws.Range("W113").Value = 1
$ws.Range("X113").Value = 1

$ws.Range("R114").Value = 50
$ws.Range("S114").Value = 1800
$ws.Range("W114").Value = 1
$ws.Range("X114").Value = 1

$ws.Range("R115").Value = 50
$ws.Range("S115").Value = 1800
$ws.Range("W115").Value = 1
$ws.Range("X115").Value = 1

$ws.Range("R116").Value = 50
$ws.Range("S116").Value = 2100
$ws.Range("W116").Value = 1
$ws.Range("X116").Value = 1

$ws.Range("R117").Value = 50
$ws.Range("S117").Value = 2200
$ws.Range("W117").Value = 1
$ws.Range("X117").Value = 1

$ws.Range("R118").Value = 50
$ws.Range("S118").Value = 2300
$ws.Range("W118").Value = 1
$ws.Range("X118").Value = 1

$ws.Range("R119").Value = 50
$ws.Range("S119").Value = 2400
$ws.Range("W119").Value = 1
$ws.Range("X119").Value = 1

$ws.Range("R120").Value = 0
$ws.Range("S120").Value = 1700
$ws.Range("W120").Value = 1
$ws.Range("X120").Value = 1

$ws.Range("R121").Value = 0
$ws.Range("S121").Value = 1800
$ws.Range("W121").Value = 1
$ws.Range("X121").Value = 1

$ws.Range("R122").Value = 0
$ws.Range("S122").Value = 1900
$ws.Range("W122").Value = 1
$ws.Range("X122").Value = 1

$ws.Range("R123").Value = 0
$ws.Range("S123").Value = 2100
$ws.Range("W123").Value = 1
$ws.Range("X123").Value = 1

$ws.Range("R124").Value = 0
$ws.Range("S124").Value = 2200
$ws.Range("W124").Value = 1
$ws.Range("X124").Value = 1

$ws.Range("R125").Value = 0
$ws.Range("S125").Value = 2300
$ws.Range("W125").Value = 1
$ws.Range("X125").Value = 1
